# Generate Report for Handoff
# Refreshes the localization-status report: the tracked files changed from
# the old png/md fixtures to the new calleeMd1/calleeMd2/callerMd1/callerMd2
# fixtures, and a 4th file (callerMd2.md) is now tracked, adding a row to
# every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "calleeMd1.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-22 11:05:09"

$ov.Range("A3").Value = "calleeMd2.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-22 11:05:09"

$ov.Range("A4").Value = "callerMd1.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-22 11:05:09"

$ov.Range("A5").Value = "callerMd2.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-03-22 11:05:09"
$ov.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild the hyperlinks on column A in row order (Delete() on the
# collection clears every link on the sheet, so we just re-add them all).
$ov.Range("A1").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd1.md", "", "", "callerMd1.md")
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd2.md", "", "", "callerMd2.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "calleeMd1.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-22 11:05:05"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"
$zh.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$zh.Range("A3").Value = "calleeMd2.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-22 11:05:05"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"
$zh.Range("K3").Value = "e2e\callerMd1.md"

$zh.Range("A4").Value = "callerMd1.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-22 11:05:05"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$zh.Range("J4").Value = "Include"

$zh.Range("A5").Value = "callerMd2.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf"
$zh.Range("E5").Value = "2016-03-22 11:05:05"
$zh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H5").Value = "0001-01-01 00:00:00"
$zh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("I5").Value = "e2e\calleeMd1.md"
$zh.Range("J5").Value = "Include"

$zh.Range("A1").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd1.md", "", "", "callerMd1.md")
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd2.md", "", "", "callerMd2.md")
$zh.Hyperlinks.Add($zh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "calleeMd1.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("D2").Value = "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf"
$de.Range("E2").Value = "2016-03-22 11:05:09"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"
$de.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$de.Range("A3").Value = "calleeMd2.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf"
$de.Range("E3").Value = "2016-03-22 11:05:09"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"
$de.Range("K3").Value = "e2e\callerMd1.md"

$de.Range("A4").Value = "callerMd1.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf"
$de.Range("E4").Value = "2016-03-22 11:05:09"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$de.Range("J4").Value = "Include"

$de.Range("A5").Value = "callerMd2.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf"
$de.Range("E5").Value = "2016-03-22 11:05:09"
$de.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H5").Value = "0001-01-01 00:00:00"
$de.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I5").Value = "e2e\calleeMd1.md"
$de.Range("J5").Value = "Include"

$de.Range("A1").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd1.md", "", "", "callerMd1.md")
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e/callerMd2.md", "", "", "callerMd2.md")
$de.Hyperlinks.Add($de.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf")
